$p = $ppt.ActivePresentation

# --- Slide 1 (title slide): the title and subtitle runs already hold
# the correct text; PowerPoint's real "re-type" edit collapses the
# redundant trailing <a:endParaRPr> once the paragraph's sole run
# carries the paragraph-mark formatting. Reproduce that by clearing
# and re-inserting each run so the stray endParaRPr is dropped.
$s1 = $p.Slides.Item(1)

$titleRange = $s1.Shapes.Item(1).TextFrame.TextRange
$titleRange.Delete()
$null = $titleRange.InsertAfter("Pattern Extraction")

$subtitleRange = $s1.Shapes.Item(2).TextFrame.TextRange
$subtitleRange.Delete()
$null = $subtitleRange.InsertAfter("Stauffer Guy-Raphaël and Chevalley Gibran")

# --- Slide 9 (last slide, "Sources"): title placeholder was left
# empty; give it its missing title text.
$s9 = $p.Slides.Item(9)
$sourcesTitle = $s9.Shapes.Item(1).TextFrame.TextRange
$sourcesRun = $sourcesTitle.InsertAfter("Sources")
$sourcesRun.LanguageID = "fr-CH"
